$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historico")

$row = 34

$ws.Cells.Item($row, 1).Value = "04/01/2026 23:11:46"
$ws.Cells.Item($row, 2).Value = "04/01 23:00"
$ws.Cells.Item($row, 3).Value = "Folha de S.Paulo - Poder - Principal"
$ws.Cells.Item($row, 4).Value = "Governo Lula chega a 2026 com medidas eleitorais pendentes no Congresso"
$ws.Cells.Item($row, 5).Value = "https://redir.folha.com.br/redir/online/poder/rss091/*https://www1.folha.uol.com.br/poder/2026/01/governo-lula-chega-a-2026-com-medidas-eleitorais-pendentes-no-congresso.shtml"
$ws.Cells.Item($row, 6).Value = "congresso"
$ws.Cells.Item($row, 7).Value = "/folha-topicos/pt/`"&gt;PT&lt;/a&gt;) chega a 2026 com medidas do seu pacote eleitoral pendentes no Congresso após um &lt;a href=`"https://www1.folha.uol.com.br/poder/2025/12/motta-encerra-1o-ano-a-frent"
